# Cluster Keywords.xlsx - "Add files via upload"
#
# The author added two new keyword rows to the Cluster_Keywords table:
#   - "Kids"  (Length 4) under the existing "Childcare" category, inserted
#     as the new sheet row 7 (pushing the remaining Childcare/… rows down).
#   - "Polic" (Length 5) under a brand-new "Government" category, inserted
#     right after "Serco"/Office and before "Stati"/Public Transport.
#
# Both insertions shift every row below them down by one, which is why the
# table/autofilter/dimension grow from A1:C73 to A1:C75 and the
# conditional-formatting ranges shift accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert "Kids" / "Childcare" as the new row 7
# ---------------------------------------------------------------------
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).Value = "Kids"
$ws.Cells.Item(7, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Cells.Item(7, 3).Value = "Childcare"

# Copy the formatting (style) from the row above so the new row matches
# the rest of the table's look (style index 1 in the original file).
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Insert "Polic" / "Government" as the new row 69
#    (after the first insertion, old row 68 "Serco"/Office now sits at
#    row 68, and old row 69 "Stati"/Public Transport now sits at row 69 -
#    the new row goes in before that, i.e. at row 69).
# ---------------------------------------------------------------------
$ws.Rows.Item(69).Insert()
$ws.Cells.Item(69, 1).Value = "Polic"
$ws.Cells.Item(69, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Cells.Item(69, 3).Value = "Government"

$ws.Range("A68:C68").Copy()
$ws.Range("A69:C69").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Grow the table (and its autofilter) to cover the two new rows
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C75"))

# Resize can regenerate the calculated column with an alternate (but
# equivalent) structured-reference syntax on the newly covered rows;
# normalize every Length cell back to the original formula text so the
# whole column is consistent again.
for ($r = 2; $r -le 75; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
}

# ---------------------------------------------------------------------
# 4) Fix up the conditional-formatting ranges that referenced fixed rows
#    below the insertion points, so they keep pointing at the same
#    logical cells (Excel does this automatically on a real row insert).
# ---------------------------------------------------------------------
$cf = $ws.Cells.FormatConditions

$cf.Item(1).ModifyAppliesToRange($ws.Range("A31:A32"))
$cf.Item(2).ModifyAppliesToRange($ws.Range("A24:A25"))
$cf.Item(3).ModifyAppliesToRange($ws.Range("B2:B75"))

$r1 = $ws.Range("A2:A23")
$r2 = $ws.Range("A62:A75")
$r3 = $ws.Range("A26:A30")
$r4 = $ws.Range("A33:A60")
$u = $excel.Union($excel.Union($r1, $r2), $excel.Union($r3, $r4))
$cf.Item(4).ModifyAppliesToRange($u)

# ---------------------------------------------------------------------
# 5) Restore the selection that was active after the edit
# ---------------------------------------------------------------------
$ws.Range("A69").Select()
$excel.ActiveWindow.ScrollRow = 43

Write-Output "Applied Kids/Childcare + Polic/Government row inserts"
